$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = 45978
$ws.Range("A29").NumberFormat = "d-mmm-yy"
$ws.Range("B29").Value = 5610
$ws.Range("C29").Value = 4222
$ws.Range("D29").Value = 3716
$ws.Range("E29").Value = 331
$ws.Range("F29").Value = 94
$ws.Range("G29").Value = 70
$ws.Range("H29").Value = 11
$ws.Range("I29").Value = 0

$ws.Range("A29:I29").Select()
